# Updates cryptocurrency price (D) and 1h-volume-change (E) columns
# with freshly scraped figures. Numeric-looking price strings are
# entered with a leading quote-prefix ($Quote) so Excel stores them
# as text (matching the original inline-string cells) instead of
# silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$Quote = "'"

$ws.Range('D2').Value = '63.515.52'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '2.651.63'
$ws.Range('E3').Value = '  +2.45%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = $Quote + '591.06'
$ws.Range('E5').Value = '  +1.16%  '
$ws.Range('D6').Value = $Quote + '144.15'
$ws.Range('E6').Value = '  -1.79%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = $Quote + '0.588'
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('D9').Value = '2.650.63'
$ws.Range('E9').Value = '  +2.47%  '
$ws.Range('D10').Value = $Quote + '0.107'
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('D11').Value = $Quote + '5.61'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').Value = $Quote + '0.353'
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('D14').Value = $Quote + '27.40'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').Value = '3.125.70'
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range('D16').Value = '63.422.24'
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('D18').Value = '2.611.35'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = $Quote + '11.37'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').Value = $Quote + '340.41'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').Value = $Quote + '4.35'
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('D22').Value = $Quote + '6.71'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = $Quote + '67.55'
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').Value = $Quote + '1.63'
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('E26').Value = '  +6.62%  '
$ws.Range('D27').Value = $Quote + '0.166'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').Value = $Quote + '550.90'
$ws.Range('E28').Value = '  +16.77%  '
$ws.Range('D29').Value = $Quote + '0.999'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = $Quote + '8.40'
$ws.Range('E30').Value = '  +1.09%  '
$ws.Range('D31').Value = $Quote + '7.75'
$ws.Range('E31').Value = '  -1.24%  '
$ws.Range('D32').Value = $Quote + '1.81'
$ws.Range('E32').Value = '  +12.95%  '
$ws.Range('E33').Value = '  +1.99%  '
$ws.Range('D34').Value = '0.0₃0806'
$ws.Range('E34').Value = '  -1.96%  '
$ws.Range('D35').Value = $Quote + '173.94'
$ws.Range('E35').Value = '  -1.70%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +6.77%  '
$ws.Range('D38').Value = $Quote + '0.402'
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('D39').Value = $Quote + '19.07'
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('D40').Value = $Quote + '1.79'
$ws.Range('E40').Value = '  +5.01%  '
$ws.Range('D41').Value = $Quote + '171.05'
$ws.Range('E41').Value = '  +7.82%  '
$ws.Range('D42').Value = $Quote + '1.00'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = $Quote + '40.24'
$ws.Range('E43').Value = '  +1.67%  '
$ws.Range('D44').Value = $Quote + '3.73'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').Value = $Quote + '22.10'
$ws.Range('E45').Value = '  +4.43%  '
$ws.Range('D46').Value = $Quote + '0.627'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('D47').Value = $Quote + '0.0551'
$ws.Range('E47').Value = '  +1.69%  '
$ws.Range('D48').Value = $Quote + '0.0959'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('E50').Value = '  +1.34%  '
$ws.Range('E51').Value = '  -0.79%  '
